$d = $word.ActiveDocument

$d.Content.Find.Execute("40÷8=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=17, 2", 2) | Out-Null
$d.Content.Find.Execute("30÷3=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷9=8, 6", 2) | Out-Null
$d.Content.Find.Execute("89÷8=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=10, 7", 2) | Out-Null
$d.Content.Find.Execute("69÷5=13, 4", $true, $false, $false, $false, $false, $true, 1, $false, "64÷8=8, 0", 2) | Out-Null
$d.Content.Find.Execute("85÷9=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "44÷7=6, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷7=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 2) | Out-Null
$d.Content.Find.Execute("49÷7=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷2=9, 0", 2) | Out-Null
$d.Content.Find.Execute("70÷2=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=38, 0", 2) | Out-Null
$d.Content.Find.Execute("61÷5=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷7=7, 3", 2) | Out-Null
$d.Content.Find.Execute("69÷4=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2) | Out-Null
$d.Content.Find.Execute("67÷5=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=17, 1", 2) | Out-Null
$d.Content.Find.Execute("63÷3=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "86÷6=14, 2", 2) | Out-Null
$d.Content.Find.Execute("83÷7=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "79÷6=13, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷2=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "77÷3=25, 2", 2) | Out-Null
$d.Content.Find.Execute("52÷5=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "73÷7=10, 3", 2) | Out-Null
$d.Content.Find.Execute("12÷7=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "92÷7=13, 1", 2) | Out-Null
$d.Content.Find.Execute("39÷2=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=18, 2", 2) | Out-Null
$d.Content.Find.Execute("58÷3=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷6=2, 2", 2) | Out-Null
$d.Content.Find.Execute("59÷3=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "47÷4=11, 3", 2) | Out-Null
$d.Content.Find.Execute("25÷6=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "31÷9=3, 4", 2) | Out-Null
$d.Content.Find.Execute("55÷8=6, 7", $true, $false, $false, $false, $false, $true, 1, $false, "62÷7=8, 6", 2) | Out-Null
$d.Content.Find.Execute("38÷5=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=7, 6", 2) | Out-Null
$d.Content.Find.Execute("43÷8=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "88÷6=14, 4", 2) | Out-Null
$d.Content.Find.Execute("63÷5=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "66÷6=11, 0", 2) | Out-Null
$d.Content.Find.Execute("81÷2=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "48÷3=16, 0", 2) | Out-Null
